$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 2967  # H33: 4350.5 -> 2967
$ws.Cells.Item(33, 9).Value = 2967  # I33: 4350.5 -> 2967
$ws.Cells.Item(33, 11).Value = 2967  # K33: 4350.5 -> 2967
$ws.Cells.Item(33, 13).Value = -2738  # M33: -4121.5 -> -2738
$ws.Cells.Item(51, 8).Value = 5999.6665  # H51: 5685.4287 -> 5999.6665
$ws.Cells.Item(51, 9).Value = 4166.6665  # I51: 4075 -> 4166.6665
$ws.Cells.Item(51, 11).Value = 4166.6665  # K51: 4075 -> 4166.6665
$ws.Cells.Item(51, 13).Value = -3682.6665  # M51: -3591 -> -3682.6665
$ws.Cells.Item(76, 8).Value = 3550.6  # H76: 3590.75 -> 3550.6
$ws.Cells.Item(76, 10).Value = 3390  # J76: 0 -> 3390
$ws.Cells.Item(76, 12).Value = 3390  # L76: 0 -> 3390
$ws.Cells.Item(76, 14).Value = -4020  # N76: None -> -4020
$ws.Cells.Item(79, 8).Value = 3550.6  # H79: 3590.75 -> 3550.6
$ws.Cells.Item(79, 10).Value = 3390  # J79: 0 -> 3390
$ws.Cells.Item(79, 12).Value = 3390  # L79: 0 -> 3390
$ws.Cells.Item(79, 14).Value = -5574  # N79: None -> -5574
$ws.Cells.Item(80, 8).Value = 8730.333000000001  # H80: 7084.5557 -> 8730.333000000001
$ws.Cells.Item(80, 9).Value = 10251  # I80: 6975.25 -> 10251
$ws.Cells.Item(80, 10).Value = 7970  # J80: 7172 -> 7970
$ws.Cells.Item(80, 11).Value = 30753  # K80: 20925.75 -> 30753
$ws.Cells.Item(80, 12).Value = 23910  # L80: 21516 -> 23910
$ws.Cells.Item(80, 13).Value = -29755  # M80: -19927.75 -> -29755
$ws.Cells.Item(80, 14).Value = -25906  # N80: -23512 -> -25906
$ws.Cells.Item(83, 8).Value = 8730.333000000001  # H83: 7084.5557 -> 8730.333000000001
$ws.Cells.Item(83, 9).Value = 10251  # I83: 6975.25 -> 10251
$ws.Cells.Item(83, 10).Value = 7970  # J83: 7172 -> 7970
$ws.Cells.Item(83, 11).Value = 92259  # K83: 62777.25 -> 92259
$ws.Cells.Item(83, 12).Value = 71730  # L83: 64548 -> 71730
$ws.Cells.Item(83, 13).Value = -87267  # M83: -57785.25 -> -87267
$ws.Cells.Item(83, 14).Value = -81714  # N83: -74532 -> -81714
$ws.Cells.Item(107, 8).Value = 703.375  # H107: 701.75 -> 703.375
$ws.Cells.Item(107, 9).Value = 565.8  # I107: 532.5 -> 565.8
$ws.Cells.Item(107, 10).Value = 932.6667  # J107: 871 -> 932.6667
$ws.Cells.Item(107, 11).Value = 565.8  # K107: 532.5 -> 565.8
$ws.Cells.Item(107, 12).Value = 932.6667  # L107: 871 -> 932.6667
$ws.Cells.Item(107, 13).Value = 1354.2  # M107: 1387.5 -> 1354.2
$ws.Cells.Item(107, 14).Value = -4772.6667  # N107: -4711 -> -4772.6667
$ws.Cells.Item(127, 8).Value = 3228  # H127: 5480 -> 3228
$ws.Cells.Item(127, 9).Value = 3228  # I127: 5480 -> 3228
$ws.Cells.Item(127, 11).Value = 9684  # K127: 16440 -> 9684
$ws.Cells.Item(127, 13).Value = -4724  # M127: -11480 -> -4724
$ws.Cells.Item(129, 8).Value = 2700.158  # H129: 3130.3684 -> 2700.158
$ws.Cells.Item(129, 9).Value = 783.2  # I129: 2331.6667 -> 783.2
$ws.Cells.Item(129, 10).Value = 3384.7856  # J129: 3499 -> 3384.7856
$ws.Cells.Item(129, 11).Value = 2349.6  # K129: 6995.000100000001 -> 2349.6
$ws.Cells.Item(129, 12).Value = 10154.3568  # L129: 10497 -> 10154.3568
$ws.Cells.Item(129, 13).Value = 2650.4  # M129: -1995.000100000001 -> 2650.4
$ws.Cells.Item(129, 14).Value = -20154.3568  # N129: -20497 -> -20154.3568
$ws.Cells.Item(131, 8).Value = 31448.5  # H131: 18273.285 -> 31448.5
$ws.Cells.Item(131, 9).Value = 50899.5  # I131: 20783.6 -> 50899.5
$ws.Cells.Item(131, 11).Value = 152698.5  # K131: 62350.8 -> 152698.5
$ws.Cells.Item(131, 13).Value = -147658.5  # M131: -57310.8 -> -147658.5
$ws.Cells.Item(141, 8).Value = 3184.5386  # H141: 2866.5833 -> 3184.5386
$ws.Cells.Item(141, 10).Value = 6500  # J141: 6000 -> 6500
$ws.Cells.Item(141, 12).Value = 19500  # L141: 18000 -> 19500
$ws.Cells.Item(141, 14).Value = -29860  # N141: -28360 -> -29860

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3023.9143  # H32: 2987.361 -> 3023.9143
$ws.Cells.Item(32, 9).Value = 3132.5173  # I32: 3085.0334 -> 3132.5173
$ws.Cells.Item(32, 11).Value = 3132.5173  # K32: 3085.0334 -> 3132.5173
$ws.Cells.Item(32, 13).Value = -2845.5173  # M32: -2798.0334 -> -2845.5173
$ws.Cells.Item(63, 8).Value = 6540.5  # H63: 6507 -> 6540.5
$ws.Cells.Item(63, 9).Value = 6088.6665  # I63: 6077.125 -> 6088.6665
$ws.Cells.Item(63, 11).Value = 6088.6665  # K63: 6077.125 -> 6088.6665
$ws.Cells.Item(63, 13).Value = -5402.6665  # M63: -5391.125 -> -5402.6665
$ws.Cells.Item(66, 8).Value = 6540.5  # H66: 6507 -> 6540.5
$ws.Cells.Item(66, 9).Value = 6088.6665  # I66: 6077.125 -> 6088.6665
$ws.Cells.Item(66, 11).Value = 30443.3325  # K66: 30385.625 -> 30443.3325
$ws.Cells.Item(66, 13).Value = -27011.3325  # M66: -26953.625 -> -27011.3325
$ws.Cells.Item(122, 8).Value = 0  # H122: 2406 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 12 -> 0
$ws.Cells.Item(122, 10).Value = 0  # J122: 4800 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 36 -> 0
$ws.Cells.Item(122, 12).Value = 0  # L122: 14400 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: 2414 -> (removed)
$ws.Cells.Item(122, 14).ClearContents()  # N122: -19300 -> (removed)
$ws.Cells.Item(132, 8).Value = 3497.6155  # H132: 3254.6155 -> 3497.6155
$ws.Cells.Item(132, 9).Value = 3405  # I132: 3191.3333 -> 3405
$ws.Cells.Item(132, 10).Value = 4007  # J132: 4014 -> 4007
$ws.Cells.Item(132, 11).Value = 10215  # K132: 9573.999899999999 -> 10215
$ws.Cells.Item(132, 12).Value = 12021  # L132: 12042 -> 12021
$ws.Cells.Item(132, 13).Value = -7685  # M132: -7043.999899999999 -> -7685
$ws.Cells.Item(132, 14).Value = -17081  # N132: -17102 -> -17081

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2699.625  # H31: 1873.9642 -> 2699.625
$ws.Cells.Item(31, 9).Value = 3493.625  # I31: 2020.8889 -> 3493.625
$ws.Cells.Item(31, 10).Value = 1905.625  # J31: 1609.5 -> 1905.625
$ws.Cells.Item(31, 11).Value = 3493.625  # K31: 2020.8889 -> 3493.625
$ws.Cells.Item(31, 12).Value = 1905.625  # L31: 1609.5 -> 1905.625
$ws.Cells.Item(31, 13).Value = -3198.625  # M31: -1725.8889 -> -3198.625
$ws.Cells.Item(31, 14).Value = -2495.625  # N31: -2199.5 -> -2495.625
$ws.Cells.Item(34, 8).Value = 2699.625  # H34: 1873.9642 -> 2699.625
$ws.Cells.Item(34, 9).Value = 3493.625  # I34: 2020.8889 -> 3493.625
$ws.Cells.Item(34, 10).Value = 1905.625  # J34: 1609.5 -> 1905.625
$ws.Cells.Item(34, 11).Value = 3493.625  # K34: 2020.8889 -> 3493.625
$ws.Cells.Item(34, 12).Value = 1905.625  # L34: 1609.5 -> 1905.625
$ws.Cells.Item(34, 13).Value = -3291.625  # M34: -1818.8889 -> -3291.625
$ws.Cells.Item(34, 14).Value = -2309.625  # N34: -2013.5 -> -2309.625
$ws.Cells.Item(62, 8).Value = 3568.4285  # H62: 3642.923 -> 3568.4285
$ws.Cells.Item(62, 9).Value = 3448.5  # I62: 3569.7144 -> 3448.5
$ws.Cells.Item(62, 11).Value = 3448.5  # K62: 3569.7144 -> 3448.5
$ws.Cells.Item(62, 13).Value = -2824.5  # M62: -2945.7144 -> -2824.5
$ws.Cells.Item(65, 8).Value = 3568.4285  # H65: 3642.923 -> 3568.4285
$ws.Cells.Item(65, 9).Value = 3448.5  # I65: 3569.7144 -> 3448.5
$ws.Cells.Item(65, 11).Value = 17242.5  # K65: 17848.572 -> 17242.5
$ws.Cells.Item(65, 13).Value = -14122.5  # M65: -14728.572 -> -14122.5
$ws.Cells.Item(105, 8).Value = 2854.375  # H105: 2860.9375 -> 2854.375
$ws.Cells.Item(105, 9).Value = 1546.25  # I105: 1559.375 -> 1546.25
$ws.Cells.Item(105, 11).Value = 1546.25  # K105: 1559.375 -> 1546.25
$ws.Cells.Item(105, 13).Value = 200.75  # M105: 187.625 -> 200.75
$ws.Cells.Item(122, 8).Value = 2600.6  # H122: 2828.2144 -> 2600.6
$ws.Cells.Item(122, 9).Value = 2108.4443  # I122: 2264 -> 2108.4443
$ws.Cells.Item(122, 10).Value = 3338.8333  # J122: 3843.8 -> 3338.8333
$ws.Cells.Item(122, 11).Value = 6325.3329  # K122: 6792 -> 6325.3329
$ws.Cells.Item(122, 12).Value = 10016.4999  # L122: 11531.4 -> 10016.4999
$ws.Cells.Item(122, 13).Value = -3875.3329  # M122: -4342 -> -3875.3329
$ws.Cells.Item(122, 14).Value = -14916.4999  # N122: -16431.4 -> -14916.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 377.64285  # H122: 361.3125 -> 377.64285
$ws.Cells.Item(122, 9).Value = 225.42857  # I122: 230.22223 -> 225.42857
$ws.Cells.Item(122, 11).Value = 2028.85713  # K122: 2072.00007 -> 2028.85713
$ws.Cells.Item(122, 13).Value = 421.1428699999999  # M122: 377.9999299999999 -> 421.1428699999999
$ws.Cells.Item(129, 8).Value = 558160.75  # H129: 590977.25 -> 558160.75
$ws.Cells.Item(129, 9).Value = 1378.875  # I129: 1535.7142 -> 1378.875
$ws.Cells.Item(129, 11).Value = 4136.625  # K129: 4607.142599999999 -> 4136.625
$ws.Cells.Item(129, 13).Value = 863.375  # M129: 392.8574000000008 -> 863.375
$ws.Cells.Item(139, 8).Value = 727.9231  # H139: 753.3333 -> 727.9231
$ws.Cells.Item(139, 9).Value = 607.44446  # I139: 630.5 -> 607.44446
$ws.Cells.Item(139, 11).Value = 1822.33338  # K139: 1891.5 -> 1822.33338
$ws.Cells.Item(139, 13).Value = 3317.66662  # M139: 3248.5 -> 3317.66662
$ws.Cells.Item(140, 8).Value = 10340.177  # H140: 10924.0625 -> 10340.177
$ws.Cells.Item(140, 9).Value = 1118  # I140: 1148 -> 1118
$ws.Cells.Item(140, 11).Value = 3354  # K140: 3444 -> 3354
$ws.Cells.Item(140, 13).Value = 1826  # M140: 1736 -> 1826

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6879.4  # H70: 5223.375 -> 6879.4
$ws.Cells.Item(70, 9).Value = 6599.5  # I70: 4826.857 -> 6599.5
$ws.Cells.Item(70, 11).Value = 6599.5  # K70: 4826.857 -> 6599.5
$ws.Cells.Item(70, 13).Value = -6329.5  # M70: -4556.857 -> -6329.5
$ws.Cells.Item(73, 8).Value = 6879.4  # H73: 5223.375 -> 6879.4
$ws.Cells.Item(73, 9).Value = 6599.5  # I73: 4826.857 -> 6599.5
$ws.Cells.Item(73, 11).Value = 6599.5  # K73: 4826.857 -> 6599.5
$ws.Cells.Item(73, 13).Value = -5663.5  # M73: -3890.857 -> -5663.5
$ws.Cells.Item(80, 8).Value = 4164  # H80: 3748 -> 4164
$ws.Cells.Item(80, 9).Value = 4164  # I80: 3748 -> 4164
$ws.Cells.Item(80, 11).Value = 4164  # K80: 3748 -> 4164
$ws.Cells.Item(80, 13).Value = -3166  # M80: -2750 -> -3166
$ws.Cells.Item(83, 8).Value = 4164  # H83: 3748 -> 4164
$ws.Cells.Item(83, 9).Value = 4164  # I83: 3748 -> 4164
$ws.Cells.Item(83, 11).Value = 20820  # K83: 18740 -> 20820
$ws.Cells.Item(83, 13).Value = -15828  # M83: -13748 -> -15828
$ws.Cells.Item(122, 8).Value = 3399.389  # H122: 3399.5 -> 3399.389
$ws.Cells.Item(122, 9).Value = 3472.3635  # I122: 3472.5454 -> 3472.3635
$ws.Cells.Item(122, 11).Value = 10417.0905  # K122: 10417.6362 -> 10417.0905
$ws.Cells.Item(122, 13).Value = -7967.0905  # M122: -7967.636200000001 -> -7967.0905
$ws.Cells.Item(126, 8).Value = 4499.3335  # H126: 3799.6667 -> 4499.3335
$ws.Cells.Item(126, 9).Value = 4499.3335  # I126: 3799.6667 -> 4499.3335
$ws.Cells.Item(126, 11).Value = 13498.0005  # K126: 11399.0001 -> 13498.0005
$ws.Cells.Item(126, 13).Value = -11028.0005  # M126: -8929.000100000001 -> -11028.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1648.75  # H22: 1969.8 -> 1648.75
$ws.Cells.Item(22, 9).Value = 1347.5  # I22: 1949.6666 -> 1347.5
$ws.Cells.Item(22, 10).Value = 1950  # J22: 2000 -> 1950
$ws.Cells.Item(22, 11).Value = 1347.5  # K22: 1949.6666 -> 1347.5
$ws.Cells.Item(22, 12).Value = 1950  # L22: 2000 -> 1950
$ws.Cells.Item(22, 13).Value = -1052.5  # M22: -1654.6666 -> -1052.5
$ws.Cells.Item(22, 14).Value = -2540  # N22: -2590 -> -2540
$ws.Cells.Item(27, 8).Value = 1648.75  # H27: 1969.8 -> 1648.75
$ws.Cells.Item(27, 9).Value = 1347.5  # I27: 1949.6666 -> 1347.5
$ws.Cells.Item(27, 10).Value = 1950  # J27: 2000 -> 1950
$ws.Cells.Item(27, 11).Value = 1347.5  # K27: 1949.6666 -> 1347.5
$ws.Cells.Item(27, 12).Value = 1950  # L27: 2000 -> 1950
$ws.Cells.Item(27, 13).Value = -1240.5  # M27: -1842.6666 -> -1240.5
$ws.Cells.Item(27, 14).Value = -2164  # N27: -2214 -> -2164
$ws.Cells.Item(40, 8).Value = 3160.7693  # H40: 3682.3333 -> 3160.7693
$ws.Cells.Item(40, 9).Value = 2826.3635  # I40: 3254.5557 -> 2826.3635
$ws.Cells.Item(40, 10).Value = 5000  # J40: 4965.6665 -> 5000
$ws.Cells.Item(40, 11).Value = 2826.3635  # K40: 3254.5557 -> 2826.3635
$ws.Cells.Item(40, 12).Value = 5000  # L40: 4965.6665 -> 5000
$ws.Cells.Item(40, 13).Value = -2690.3635  # M40: -3118.5557 -> -2690.3635
$ws.Cells.Item(40, 14).Value = -5272  # N40: -5237.6665 -> -5272
$ws.Cells.Item(55, 8).Value = 558.2  # H55: 199 -> 558.2
$ws.Cells.Item(55, 9).Value = 197.25  # I55: 199 -> 197.25
$ws.Cells.Item(55, 10).Value = 2002  # J55: 0 -> 2002
$ws.Cells.Item(55, 11).Value = 197.25  # K55: 199 -> 197.25
$ws.Cells.Item(55, 12).Value = 2002  # L55: 0 -> 2002
$ws.Cells.Item(55, 13).Value = -24.25  # M55: -26 -> -24.25
$ws.Cells.Item(55, 14).Value = -2348  # N55: None -> -2348
$ws.Cells.Item(132, 8).Value = 4827.8125  # H132: 5403.143 -> 4827.8125
$ws.Cells.Item(132, 9).Value = 5594.5835  # I132: 5637 -> 5594.5835
$ws.Cells.Item(132, 10).Value = 2527.5  # J132: 4000 -> 2527.5
$ws.Cells.Item(132, 11).Value = 16783.7505  # K132: 16911 -> 16783.7505
$ws.Cells.Item(132, 12).Value = 7582.5  # L132: 12000 -> 7582.5
$ws.Cells.Item(132, 13).Value = -14253.7505  # M132: -14381 -> -14253.7505
$ws.Cells.Item(132, 14).Value = -12642.5  # N132: -17060 -> -12642.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2501350  # H81: 1429678.4 -> 2501350
$ws.Cells.Item(81, 9).Value = 1199  # I81: 949.4 -> 1199
$ws.Cells.Item(81, 11).Value = 2398  # K81: 1898.8 -> 2398
$ws.Cells.Item(81, 13).Value = -1337  # M81: -837.8 -> -1337
$ws.Cells.Item(84, 8).Value = 2501350  # H84: 1429678.4 -> 2501350
$ws.Cells.Item(84, 9).Value = 1199  # I84: 949.4 -> 1199
$ws.Cells.Item(84, 11).Value = 11990  # K84: 9494 -> 11990
$ws.Cells.Item(84, 13).Value = -6686  # M84: -4190 -> -6686
